# Auto-generated edit script for violent-crime-ytd workbook
# Commit: Add data for 2024-10-18
# Updates the 2024 (column K) values across the Citywide Totals, By Neighborhood,
# and individual neighborhood sheets to reflect the latest data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6310
$ws.Range("K3").Value = 6511
$ws.Range("K6").Value = 7186
$ws.Range("K7").Value = 21826

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 190
$ws.Range("K6").Value = 150
$ws.Range("K7").Value = 641
$ws.Range("K10").Value = 127
$ws.Range("K11").Value = 406
$ws.Range("K14").Value = 110
$ws.Range("K18").Value = 145
$ws.Range("K19").Value = 638
$ws.Range("K20").Value = 520
$ws.Range("K21").Value = 70
$ws.Range("K24").Value = 67
$ws.Range("K25").Value = 106
$ws.Range("K27").Value = 208
$ws.Range("K29").Value = 1184
$ws.Range("K31").Value = 242
$ws.Range("K33").Value = 953
$ws.Range("K34").Value = 125
$ws.Range("K36").Value = 278
$ws.Range("K37").Value = 741
$ws.Range("K41").Value = 152
$ws.Range("K42").Value = 805
$ws.Range("K45").Value = 31
$ws.Range("K51").Value = 279
$ws.Range("K52").Value = 574
$ws.Range("K53").Value = 279
$ws.Range("K54").Value = 430
$ws.Range("K59").Value = 39
$ws.Range("K60").Value = 128
$ws.Range("K63").Value = 61
$ws.Range("K65").Value = 512
$ws.Range("K67").Value = 854
$ws.Range("K70").Value = 38
$ws.Range("K73").Value = 197
$ws.Range("K77").Value = 151
$ws.Range("K78").Value = 247
$ws.Range("K79").Value = 552
$ws.Range("K83").Value = 473
$ws.Range("K84").Value = 175
$ws.Range("K85").Value = 1015
$ws.Range("K88").Value = 232
$ws.Range("K89").Value = 322
$ws.Range("K90").Value = 205
$ws.Range("K91").Value = 252
$ws.Range("K94").Value = 292
$ws.Range("K95").Value = 359
$ws.Range("K96").Value = 229
$ws.Range("K98").Value = 109
$ws.Range("K101").Value = 21826

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 37
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 209
$ws.Range("K6").Value = 177
$ws.Range("K7").Value = 641

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 105
$ws.Range("K6").Value = 132
$ws.Range("K7").Value = 406

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 92
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 322

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 332
$ws.Range("K3").Value = 348
$ws.Range("K7").Value = 1015

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 208
$ws.Range("K7").Value = 574

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 73
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 390
$ws.Range("K6").Value = 479

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 168
$ws.Range("K7").Value = 473

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 246
$ws.Range("K3").Value = 340
$ws.Range("K6").Value = 293
$ws.Range("K7").Value = 953

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 128
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 359

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 208
$ws.Range("K6").Value = 222
$ws.Range("K7").Value = 741

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 172
$ws.Range("K7").Value = 512

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 242

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 311
$ws.Range("K6").Value = 242
$ws.Range("K7").Value = 854

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 232
$ws.Range("K7").Value = 430

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 423
$ws.Range("K6").Value = 341
$ws.Range("K7").Value = 1184

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 208
$ws.Range("K7").Value = 638

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K3").Value = 31
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 243
$ws.Range("K6").Value = 302
$ws.Range("K7").Value = 805

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 74
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 247

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 121
$ws.Range("K7").Value = 252

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 184
$ws.Range("K3").Value = 178
$ws.Range("K7").Value = 552

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 179
$ws.Range("K7").Value = 520

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 83
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 292

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 190

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 70
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 76
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 31
